$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.193.87'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '2.517.62'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '590.08'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E6').Value = '  +3.74%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('E9').Value = '  +3.55%  '
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.95'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.83'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D15').Value = '67.919.98'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '2.491.10'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.04'
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.55'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '352.58'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.05'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.74'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.39%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.31'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.16'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').Value = '2.641.25'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.997'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').Value = '0.0₃0921'
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '510.66'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.90'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.122'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.48%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '165.01'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.65%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.44'
$ws.Range('D37').Style = "Normal"
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.91'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('E44').Value = '  +5.30%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '147.16'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.55'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.99%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.522'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('E48').Value = '  +3.80%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0745'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.60'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.588'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.95%  '
